$d = $word.ActiveDocument

# 1. Update activation date
$d.Content.Find.Execute(
    "Ativação: 01/01/1996", $false, $true, $false, $false, $false,
    $true, 1, $false, "Ativação: 01/01/2021", 2)

# 2. Rewrite "Objetivos" paragraph (fix leading hyphen spacing, drop trailing sentence)
$d.Content.Find.Execute(
    "-Apresentar a evolução das condições geológicas da Terra que culminaram com os recursos naturais existentes hoje, com ênfase nas reservas de combustíveis fósseis, hídricos e de minérios e a conseqüente reserva de energia advinda dessas fontes naturais. A América do Sul e do Brasil, mais especificamente, são destacados.- Acompanhar, através da literatura atual, a escassez dos recursos em vista do mau uso e das alternativas para o aproveitamento racional dos recursos existentes, com alternativas na área de geração de energia.",
    $false, $true, $false, $false, $false,
    $true, 1, $false,
    "- Apresentar a evolução das condições geológicas da Terra que culminaram com os recursos naturais existentes hoje, com ênfase nas reservas de combustíveis fósseis, hídricos e de minérios e a conseqüente reserva de energia advinda dessas fontes naturais. A América do Sul e do Brasil, mais especificamente, são destacados.",
    2)

# 3. Rewrite "Programa resumido" paragraph
$d.Content.Find.Execute(
    "- Universo e Terra,- Desenvolvimento da Terra.- Combustíveis fósseis.- Recursos renováveis e biomassa.- Recursos minerais.- Matérias-primas da grande indústria química:metais ferrosos e não-ferrosos.- Recursos hídricos.",
    $false, $true, $false, $false, $false,
    $true, 1, $false,
    "- Desenvolvimento da Terra.- Recursos minerais.- Matérias-primas da grande indústria metalúrgica: metais ferrosos e não-ferrosos",
    2)

# 4. Rewrite "Programa" paragraph
$d.Content.Find.Execute(
    "- Formação do Universo. - Formação do Sistema Solar.- Desenvolvimento da Terra.- Principais Eras Geológicas.- Petróleo.- Carvão e Gás Natural - Geração de  energia termelétrica.- Recursos Renováveis.- Biomassa ? Fontes alternativas de energia.- Matérias-primas para a grande indústria química.- Metais ferrosos.- Metais não-ferrosos.- Recursos hídricos ? Bacias hídricas.- Poluição das águas.- Escassez e reaproveitamento das águas.- Geração de energia elétrica.",
    $false, $true, $false, $false, $false,
    $true, 1, $false,
    "- Desenvolvimento da Terra. - Principais Eras Geológicas. - Matérias-primas para a grande indústria metalúrgica: metais ferrosos e metais não-ferrosos.",
    2)

# 5. Rewrite "Bibliografia" paragraph
$d.Content.Find.Execute(
    "- Schäfer, A., Fundamentos de Ecologia e Biogeografia de Águas Continentais, Ed. Universidade, Porto Alegre.- Abreu, S.F. Recursos Minerais do Brasil, Ed. Edgard Bluecher, 1973.- Carioca, J. O. B. and Arora, H.L. Biomassa-Fundamentos e Aplicações Tecnológicas, Universidade Federal do Ceará, 1984.- Fernandes, F.R.C. Quem é quem no Subsolo Brasileiro, MCT/CNPq, 1987- Petri, S. e Fúlfaro, V. Geologia do Brasil, EDUSP, 1983.- Revistas especializadas e fontes de informação multimídia as mais diversas, dado ao caráter dinâmico das informações sobre reservas minerais, geração de energia e recursos naturais em geral.",
    $false, $true, $false, $false, $false,
    $true, 1, $false,
    "- MILLER, Jr. G. T. “Ciência Ambiental”,  Editora: Thomson, 2006.- ABREU, S. F. “Recursos Minerais do Brasil”, Editora: Edgard Blücher, 1973.-  SKINNER, B. J. “Recursos Minerais da Terra”, Editora: Edgard Blücher, 1996.- WICANDER, R.; MONROE, J. S. “Fundamentos de Geologia”, Editora: Cengage Learning, 2009. - PRESS, F.; Siever, R.; Jordan, T.; Grotzinger, J. “Para Entender a Terra”, Editora: Bookman,  2006.- SCHÄFER, A. “Fundamentos de Ecologia e Biogeografia de Águas Continentais”, Editora: Universidade, Porto Alegre. - Revistas especializadas e sites, dado ao caráter dinâmico das informações sobre reservas minerais e recursos naturais em geral.",
    2)
